$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "executed at" timestamps for rows 2-54 (column B), matching the
# re-run of this MRF test case (commit: "Updated Existing MRF TC").
$newTimestamps = @(
  "Fri Oct 25 11:52:58 EDT 2024",
  "Fri Oct 25 11:53:10 EDT 2024",
  "Fri Oct 25 11:53:22 EDT 2024",
  "Fri Oct 25 11:53:34 EDT 2024",
  "Fri Oct 25 11:53:46 EDT 2024",
  "Fri Oct 25 11:53:58 EDT 2024",
  "Fri Oct 25 11:54:09 EDT 2024",
  "Fri Oct 25 11:54:21 EDT 2024",
  "Fri Oct 25 11:54:33 EDT 2024",
  "Fri Oct 25 11:54:45 EDT 2024",
  "Fri Oct 25 11:54:56 EDT 2024",
  "Fri Oct 25 11:55:09 EDT 2024",
  "Fri Oct 25 11:55:20 EDT 2024",
  "Fri Oct 25 11:55:33 EDT 2024",
  "Fri Oct 25 11:55:45 EDT 2024",
  "Fri Oct 25 11:55:58 EDT 2024",
  "Fri Oct 25 11:56:10 EDT 2024",
  "Fri Oct 25 11:56:22 EDT 2024",
  "Fri Oct 25 11:56:35 EDT 2024",
  "Fri Oct 25 11:56:47 EDT 2024",
  "Fri Oct 25 11:56:59 EDT 2024",
  "Fri Oct 25 11:57:12 EDT 2024",
  "Fri Oct 25 11:57:25 EDT 2024",
  "Fri Oct 25 11:57:38 EDT 2024",
  "Fri Oct 25 11:57:50 EDT 2024",
  "Fri Oct 25 11:58:02 EDT 2024",
  "Fri Oct 25 11:58:15 EDT 2024",
  "Fri Oct 25 11:58:27 EDT 2024",
  "Fri Oct 25 11:58:40 EDT 2024",
  "Fri Oct 25 11:58:52 EDT 2024",
  "Fri Oct 25 11:59:05 EDT 2024",
  "Fri Oct 25 11:59:17 EDT 2024",
  "Fri Oct 25 11:59:30 EDT 2024",
  "Fri Oct 25 11:59:42 EDT 2024",
  "Fri Oct 25 11:59:54 EDT 2024",
  "Fri Oct 25 12:00:06 EDT 2024",
  "Fri Oct 25 12:00:18 EDT 2024",
  "Fri Oct 25 12:00:30 EDT 2024",
  "Fri Oct 25 12:00:41 EDT 2024",
  "Fri Oct 25 12:00:53 EDT 2024",
  "Fri Oct 25 12:01:05 EDT 2024",
  "Fri Oct 25 12:01:17 EDT 2024",
  "Fri Oct 25 12:01:30 EDT 2024",
  "Fri Oct 25 12:01:42 EDT 2024",
  "Fri Oct 25 12:01:55 EDT 2024",
  "Fri Oct 25 12:02:08 EDT 2024",
  "Fri Oct 25 12:02:20 EDT 2024",
  "Fri Oct 25 12:02:33 EDT 2024",
  "Fri Oct 25 12:02:45 EDT 2024",
  "Fri Oct 25 12:02:58 EDT 2024",
  "Fri Oct 25 12:03:10 EDT 2024",
  "Fri Oct 25 12:03:23 EDT 2024",
  "Fri Oct 25 12:03:35 EDT 2024"
)

for ($i = 0; $i -lt $newTimestamps.Length; $i++) {
  $row = 2 + $i
  $ws.Cells.Item($row, 2).Value = $newTimestamps[$i]
}

